$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Title paragraph: "Functioneel ontwerp inelveropdracht 1:" ->
#    "Functioneel ontwerp inelveropdracht 2:" split across two runs
#    (a run for "Functioneel ontwerp inelveropdracht " and a run for "2").
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$oldDigit = "1"
$newDigit = "2"
$digitStart = $titlePara.Range.Start + ("Functioneel ontwerp inelveropdracht ").Length
$digitEnd = $digitStart + $oldDigit.Length
$digitRng = $d.Range($digitStart, $digitEnd)
$digitRng.Text = $newDigit
$digitRng = $d.Range($digitStart, $digitStart + $newDigit.Length)
# Force a genuine run split at this boundary (engine coalesces runs whose
# final formatting matches) by toggling Bold off/on, landing on the same
# Bold=True value the run already needs.
$digitRng.Font.Bold = $false
$digitRng.Font.Bold = $true

# ---------------------------------------------------------------------
# 2) "De speler kan vakjes aan klikken" ->
#    "De speler kan " + "in een 5x5 vakje een blokje bewegen" (two runs)
# ---------------------------------------------------------------------
$prefix = "De speler kan "
$oldTail = "vakjes aan klikken"
$newTail = "in een 5x5 vakje een blokje bewegen"
$found = $d.Content.Find.Execute($prefix + $oldTail)
if (-not $found) {
    throw "Could not find 'De speler kan vakjes aan klikken' paragraph"
}
$fullPara = $d.Paragraphs.Item(5)
$tailStart = $fullPara.Range.Start + $prefix.Length
$tailEnd = $fullPara.Range.Start + ($prefix + $oldTail).Length
$tailRng = $d.Range($tailStart, $tailEnd)
$tailRng.Text = $newTail
$tailRng = $d.Range($tailStart, $tailStart + $newTail.Length)
$tailRng.Font.Bold = $true
$tailRng.Font.Bold = $false

# ---------------------------------------------------------------------
# 3) "Na het klikken zie je een X of O, omstebeurt zet je er steeds 1
#     ergens neer" -> "Bestuur het blokje met de pijltjes toetsen."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Na het klikken zie je een X of O, omstebeurt zet je er steeds 1 ergens neer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Bestuur het blokje met de pijltjes toetsen.", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) "Het programma houd op met iets doen als alle vakjes vol zijn, dit
#     om dat er geen wincondition is." ->
#    "Zodra de speler de rand raakt, kan hij niet verder."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Het programma houd op met iets doen als alle vakjes vol zijn, dit om dat er geen wincondition is.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Zodra de speler de rand raakt, kan hij niet verder.", 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Last list item "Refresh de pagina als je weer wil spelen." loses its
#    list formatting (pStyle + numPr) and its text; only the bookmark
#    around the (now empty) paragraph remains, with a left indent.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$lastTextStart = $lastPara.Range.Start
$lastTextEnd = $lastPara.Range.End - 1
$lastTextRng = $d.Range($lastTextStart, $lastTextEnd)
$lastTextRng.Text = ""

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$lastPara.Range.ListFormat.RemoveNumbers()
$lastPara.Style = "Normal"
$lastPara.Range.ParagraphFormat.LeftIndent = 18
